$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O: header "MEMBER_ID" in O1, value 218630 in O2 (text-formatted number)
$ws.Range("O1").Value = "MEMBER_ID"
$ws.Range("O2").Value = 218630
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").HorizontalAlignment = -4131

# Reflect the new view/selection state captured in the saved workbook
$ws.Range("O2").Select()
